# Updates the "cryptos" sheet with refreshed Price (column D) and
# Volume(1h) (column E) figures, row by row, for rows 2-51.
#
# The Price/Volume cells in this sheet are stored as plain text (not
# numbers), e.g. "26.437.32" or "1.000" or "  -1.39%  ". Assigning such
# strings straight to Range.Value lets Excel "helpfully" reinterpret them
# as numbers (dropping trailing zeros, switching to scientific notation,
# etc.), which would corrupt the data. To keep them as exact text we
# temporarily force the cell's number format to Text ("@") before writing
# the value, then restore the cell style to "Normal" afterwards so no
# visible formatting change is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '26.437.32'
Set-TextValue $ws.Range('E2') '  -1.39%  '
Set-TextValue $ws.Range('D3') '1.841.30'
Set-TextValue $ws.Range('E3') '  -1.76%  '
Set-TextValue $ws.Range('E4') '  -0.04%  '
Set-TextValue $ws.Range('D5') '262.43'
Set-TextValue $ws.Range('E5') '  -5.42%  '
Set-TextValue $ws.Range('D6') '1.000'
Set-TextValue $ws.Range('E6') '  +0.00%  '
Set-TextValue $ws.Range('D7') '0.5193'
Set-TextValue $ws.Range('E7') '  -1.85%  '
Set-TextValue $ws.Range('D8') '0.3261'
Set-TextValue $ws.Range('E8') '  -4.67%  '
Set-TextValue $ws.Range('D9') '0.06786'
Set-TextValue $ws.Range('E9') '  -2.28%  '
Set-TextValue $ws.Range('D10') '18.78'
Set-TextValue $ws.Range('E10') '  -6.50%  '
Set-TextValue $ws.Range('D11') '0.7738'
Set-TextValue $ws.Range('E11') '  -3.89%  '
Set-TextValue $ws.Range('D12') '0.07727'
Set-TextValue $ws.Range('E12') '  -0.38%  '
Set-TextValue $ws.Range('D13') '1.833.34'
Set-TextValue $ws.Range('E13') '  -2.04%  '
Set-TextValue $ws.Range('D14') '88.10'
Set-TextValue $ws.Range('E14') '  -2.62%  '
Set-TextValue $ws.Range('D15') '5.003'
Set-TextValue $ws.Range('E15') '  -3.58%  '
Set-TextValue $ws.Range('D16') '0.9995'
Set-TextValue $ws.Range('D17') '13.94'
Set-TextValue $ws.Range('E17') '  -4.37%  '
Set-TextValue $ws.Range('D18') '1.0000'
Set-TextValue $ws.Range('E18') '  -0.01%  '
Set-TextValue $ws.Range('E19') '  -1.35%  '
Set-TextValue $ws.Range('D20') '26.484.20'
Set-TextValue $ws.Range('E20') '  -1.37%  '
Set-TextValue $ws.Range('D21') '2.072.91'
Set-TextValue $ws.Range('E21') '  -1.29%  '
Set-TextValue $ws.Range('D22') '4.602'
Set-TextValue $ws.Range('E22') '  -3.14%  '
Set-TextValue $ws.Range('D23') '9.538'
Set-TextValue $ws.Range('E23') '  -5.00%  '
Set-TextValue $ws.Range('D24') '6.002'
Set-TextValue $ws.Range('E24') '  -2.89%  '
Set-TextValue $ws.Range('D25') '144.91'
Set-TextValue $ws.Range('E25') '  -1.16%  '
Set-TextValue $ws.Range('E26') '  -8.07%  '
Set-TextValue $ws.Range('D27') '1.655'
Set-TextValue $ws.Range('E27') '  -0.45%  '
Set-TextValue $ws.Range('D28') '16.95'
Set-TextValue $ws.Range('E28') '  -2.31%  '
Set-TextValue $ws.Range('D29') '111.79'
Set-TextValue $ws.Range('D30') '4.200'
Set-TextValue $ws.Range('E30') '  -3.42%  '
Set-TextValue $ws.Range('D31') '4.137'
Set-TextValue $ws.Range('E31') '  -4.21%  '
Set-TextValue $ws.Range('E32') '  -2.33%  '
Set-TextValue $ws.Range('D33') '0.04818'
Set-TextValue $ws.Range('E33') '  -2.32%  '
Set-TextValue $ws.Range('D34') '1.133'
Set-TextValue $ws.Range('E34') '  -3.44%  '
Set-TextValue $ws.Range('D35') '0.7177'
Set-TextValue $ws.Range('E35') '  -1.84%  '
Set-TextValue $ws.Range('D36') '2.849'
Set-TextValue $ws.Range('E36') '  -1.21%  '
Set-TextValue $ws.Range('D37') '3.087'
Set-TextValue $ws.Range('D38') '0.01782'
Set-TextValue $ws.Range('E38') '  -4.08%  '
Set-TextValue $ws.Range('D39') '2.228'
Set-TextValue $ws.Range('E39') '  -4.41%  '
Set-TextValue $ws.Range('D40') '0.4849'
Set-TextValue $ws.Range('E40') '  -6.01%  '
Set-TextValue $ws.Range('D41') '112.66'
Set-TextValue $ws.Range('E41') '  -3.02%  '
Set-TextValue $ws.Range('D42') '0.9029'
Set-TextValue $ws.Range('E42') '  -5.17%  '
Set-TextValue $ws.Range('D43') '6.069'
Set-TextValue $ws.Range('E43') '  -1.62%  '
Set-TextValue $ws.Range('D44') '1.000'
Set-TextValue $ws.Range('E44') '  +0.04%  '
Set-TextValue $ws.Range('D45') '7.736'
Set-TextValue $ws.Range('E45') '  -4.83%  '
Set-TextValue $ws.Range('D46') '0.4158'
Set-TextValue $ws.Range('E46') '  -7.17%  '
Set-TextValue $ws.Range('D47') '0.05913'
Set-TextValue $ws.Range('E47') '  -0.49%  '
Set-TextValue $ws.Range('D48') '8.997'
Set-TextValue $ws.Range('E48') '  -3.80%  '
Set-TextValue $ws.Range('D49') '35.07'
Set-TextValue $ws.Range('E49') '  -3.65%  '
Set-TextValue $ws.Range('D50') '0.1220'
Set-TextValue $ws.Range('E50') '  -9.09%  '
Set-TextValue $ws.Range('D51') '0.8865'
Set-TextValue $ws.Range('E51') '  +0.00%  '
